# Auto-generated edit script: updates crime-count figures across neighborhood sheets
# per commit "Add data for 2022-08-11".
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Cells.Item(2, 5).Value = 46  # E2: 45 -> 46
$ws.Cells.Item(2, 7).Value = 62  # G2: 61 -> 62
$ws.Cells.Item(3, 5).Value = 94  # E3: 91 -> 94
$ws.Cells.Item(3, 8).Value = 78  # H3: 77 -> 78
$ws.Cells.Item(9, 3).Value = 304  # C9: 303 -> 304
$ws.Cells.Item(9, 4).Value = 277  # D9: 275 -> 277
$ws.Cells.Item(9, 5).Value = 265  # E9: 263 -> 265
$ws.Cells.Item(9, 8).Value = 271  # H9: 270 -> 271
$ws.Cells.Item(9, 9).Value = 339  # I9: 337 -> 339
$ws.Cells.Item(10, 2).Value = 758  # B10: 755 -> 758
$ws.Cells.Item(10, 3).Value = 924  # C10: 918 -> 924
$ws.Cells.Item(10, 4).Value = 1117  # D10: 1108 -> 1117
$ws.Cells.Item(10, 5).Value = 1363  # E10: 1355 -> 1363
$ws.Cells.Item(10, 6).Value = 1422  # F10: 1417 -> 1422
$ws.Cells.Item(10, 7).Value = 705  # G10: 703 -> 705
$ws.Cells.Item(10, 8).Value = 334  # H10: 331 -> 334
$ws.Cells.Item(10, 9).Value = 536  # I10: 534 -> 536
$ws.Cells.Item(11, 2).Value = 1084  # B11: 1081 -> 1084
$ws.Cells.Item(11, 3).Value = 1333  # C11: 1326 -> 1333
$ws.Cells.Item(11, 4).Value = 1545  # D11: 1534 -> 1545
$ws.Cells.Item(11, 5).Value = 1777  # E11: 1763 -> 1777
$ws.Cells.Item(11, 6).Value = 1927  # F11: 1922 -> 1927
$ws.Cells.Item(11, 7).Value = 1156  # G11: 1153 -> 1156
$ws.Cells.Item(11, 8).Value = 762  # H11: 757 -> 762
$ws.Cells.Item(11, 9).Value = 1101  # I11: 1097 -> 1101

$ws = $wb.Worksheets.Item("Chinatown")
$ws.Cells.Item(8, 4).Value = 6  # D8: 5 -> 6
$ws.Cells.Item(9, 4).Value = 8  # D9: 7 -> 8

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Cells.Item(2, 7).Value = 4  # G2: 3 -> 4
$ws.Cells.Item(6, 5).Value = 4  # E6: 3 -> 4
$ws.Cells.Item(7, 5).Value = 6  # E7: 5 -> 6

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Cells.Item(7, 9).Value = 21  # I7: 20 -> 21
$ws.Cells.Item(9, 7).Value = 82  # G9: 81 -> 82
$ws.Cells.Item(9, 9).Value = 61  # I9: 60 -> 61

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Cells.Item(8, 2).Value = 23  # B8: 22 -> 23
$ws.Cells.Item(8, 4).Value = 28  # D8: 27 -> 28
$ws.Cells.Item(8, 6).Value = 27  # F8: 26 -> 27
$ws.Cells.Item(8, 8).Value = 15  # H8: 13 -> 15
$ws.Cells.Item(9, 2).Value = 39  # B9: 38 -> 39
$ws.Cells.Item(9, 4).Value = 52  # D9: 51 -> 52
$ws.Cells.Item(9, 6).Value = 59  # F9: 58 -> 59
$ws.Cells.Item(9, 8).Value = 38  # H9: 36 -> 38

$ws = $wb.Worksheets.Item("Loop")
$ws.Cells.Item(8, 3).Value = 172  # C8: 170 -> 172
$ws.Cells.Item(8, 4).Value = 327  # D8: 324 -> 327
$ws.Cells.Item(8, 5).Value = 392  # E8: 389 -> 392
$ws.Cells.Item(8, 6).Value = 390  # F8: 388 -> 390
$ws.Cells.Item(8, 9).Value = 125  # I8: 124 -> 125
$ws.Cells.Item(9, 3).Value = 205  # C9: 203 -> 205
$ws.Cells.Item(9, 4).Value = 382  # D9: 379 -> 382
$ws.Cells.Item(9, 5).Value = 442  # E9: 439 -> 442
$ws.Cells.Item(9, 6).Value = 444  # F9: 442 -> 444
$ws.Cells.Item(9, 9).Value = 213  # I9: 212 -> 213

$ws = $wb.Worksheets.Item("Armour Square")
$ws.Cells.Item(3, 5).Value = 3  # E3: 2 -> 3
$ws.Cells.Item(7, 5).Value = 14  # E7: 13 -> 14

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Cells.Item(5, 5).Value = 14  # E5: 13 -> 14
$ws.Cells.Item(8, 7).Value = 71  # G8: 70 -> 71
$ws.Cells.Item(20, 5).Value = 6  # E20: 5 -> 6
$ws.Cells.Item(21, 4).Value = 8  # D21: 7 -> 8
$ws.Cells.Item(27, 4).Value = 17  # D27: 16 -> 17
$ws.Cells.Item(29, 5).Value = 15  # E29: 14 -> 15
$ws.Cells.Item(29, 7).Value = 8  # G29: 7 -> 8
$ws.Cells.Item(30, 5).Value = 26  # E30: 25 -> 26
$ws.Cells.Item(32, 7).Value = 82  # G32: 81 -> 82
$ws.Cells.Item(32, 9).Value = 61  # I32: 60 -> 61
$ws.Cells.Item(33, 3).Value = 2  # C33: 1 -> 2
$ws.Cells.Item(36, 2).Value = 39  # B36: 38 -> 39
$ws.Cells.Item(36, 4).Value = 52  # D36: 51 -> 52
$ws.Cells.Item(36, 6).Value = 59  # F36: 58 -> 59
$ws.Cells.Item(36, 8).Value = 38  # H36: 36 -> 38
$ws.Cells.Item(45, 3).Value = 11  # C45: 10 -> 11
$ws.Cells.Item(47, 3).Value = 50  # C47: 49 -> 50
$ws.Cells.Item(47, 5).Value = 45  # E47: 43 -> 45
$ws.Cells.Item(47, 9).Value = 32  # I47: 31 -> 32
$ws.Cells.Item(48, 5).Value = 9  # E48: 8 -> 9
$ws.Cells.Item(49, 8).Value = 5  # H49: 4 -> 5
$ws.Cells.Item(51, 4).Value = 5  # D51: 4 -> 5
$ws.Cells.Item(53, 3).Value = 205  # C53: 203 -> 205
$ws.Cells.Item(53, 4).Value = 382  # D53: 379 -> 382
$ws.Cells.Item(53, 5).Value = 442  # E53: 439 -> 442
$ws.Cells.Item(53, 6).Value = 444  # F53: 442 -> 444
$ws.Cells.Item(53, 9).Value = 213  # I53: 212 -> 213
$ws.Cells.Item(54, 3).Value = 6  # C54: 5 -> 6
$ws.Cells.Item(54, 8).Value = 4  # H54: 3 -> 4
$ws.Cells.Item(61, 2).Value = 4  # B61: 3 -> 4
$ws.Cells.Item(61, 6).Value = 30  # F61: 29 -> 30
$ws.Cells.Item(62, 5).Value = 19  # E62: 18 -> 19
$ws.Cells.Item(63, 5).Value = 9  # E63: 8 -> 9
$ws.Cells.Item(72, 8).Value = 4  # H72: 3 -> 4
$ws.Cells.Item(74, 3).Value = 23  # C74: 22 -> 23
$ws.Cells.Item(74, 5).Value = 52  # E74: 51 -> 52
$ws.Cells.Item(76, 4).Value = 38  # D76: 36 -> 38
$ws.Cells.Item(77, 4).Value = 35  # D77: 34 -> 35
$ws.Cells.Item(77, 9).Value = 54  # I77: 53 -> 54
$ws.Cells.Item(78, 2).Value = 22  # B78: 21 -> 22
$ws.Cells.Item(83, 4).Value = 3  # D83: 2 -> 3
$ws.Cells.Item(95, 6).Value = 13  # F95: 12 -> 13
$ws.Cells.Item(96, 5).Value = 14  # E96: 13 -> 14
$ws.Cells.Item(98, 2).Value = 1084  # B98: 1081 -> 1084
$ws.Cells.Item(98, 3).Value = 1333  # C98: 1326 -> 1333
$ws.Cells.Item(98, 4).Value = 1545  # D98: 1534 -> 1545
$ws.Cells.Item(98, 5).Value = 1777  # E98: 1763 -> 1777
$ws.Cells.Item(98, 6).Value = 1927  # F98: 1922 -> 1927
$ws.Cells.Item(98, 7).Value = 1156  # G98: 1153 -> 1156
$ws.Cells.Item(98, 8).Value = 762  # H98: 757 -> 762
$ws.Cells.Item(98, 9).Value = 1101  # I98: 1097 -> 1101

$ws = $wb.Worksheets.Item("Rush & Division")
$ws.Cells.Item(5, 2).Value = 19  # B5: 18 -> 19
$ws.Cells.Item(6, 2).Value = 22  # B6: 21 -> 22

$ws = $wb.Worksheets.Item("Lake View")
$ws.Cells.Item(2, 5).Value = 2  # E2: 1 -> 2
$ws.Cells.Item(6, 5).Value = 6  # E6: 5 -> 6
$ws.Cells.Item(6, 9).Value = 14  # I6: 13 -> 14
$ws.Cells.Item(7, 3).Value = 31  # C7: 30 -> 31
$ws.Cells.Item(8, 3).Value = 50  # C8: 49 -> 50
$ws.Cells.Item(8, 5).Value = 45  # E8: 43 -> 45
$ws.Cells.Item(8, 9).Value = 32  # I8: 31 -> 32

$ws = $wb.Worksheets.Item("Jefferson Park")
$ws.Cells.Item(5, 3).Value = 1  # C5: None -> 1
$ws.Cells.Item(7, 3).Value = 11  # C7: 10 -> 11

$ws = $wb.Worksheets.Item("Fuller Park")
$ws.Cells.Item(3, 5).Value = 2  # E3: 1 -> 2
$ws.Cells.Item(8, 7).Value = 2  # G8: 1 -> 2
$ws.Cells.Item(9, 5).Value = 15  # E9: 14 -> 15
$ws.Cells.Item(9, 7).Value = 8  # G9: 7 -> 8

$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Cells.Item(7, 4).Value = 29  # D7: 27 -> 29
$ws.Cells.Item(8, 4).Value = 38  # D8: 36 -> 38

$ws = $wb.Worksheets.Item("River North")
$ws.Cells.Item(6, 3).Value = 20  # C6: 19 -> 20
$ws.Cells.Item(6, 5).Value = 49  # E6: 48 -> 49
$ws.Cells.Item(7, 3).Value = 23  # C7: 22 -> 23
$ws.Cells.Item(7, 5).Value = 52  # E7: 51 -> 52

$ws = $wb.Worksheets.Item("Little Village")
$ws.Cells.Item(5, 4).Value = 2  # D5: 1 -> 2
$ws.Cells.Item(6, 4).Value = 5  # D6: 4 -> 5

$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Cells.Item(5, 5).Value = 5  # E5: 4 -> 5
$ws.Cells.Item(7, 5).Value = 14  # E7: 13 -> 14

$ws = $wb.Worksheets.Item("Edgewater")
$ws.Cells.Item(4, 4).Value = 5  # D4: 4 -> 5
$ws.Cells.Item(6, 4).Value = 17  # D6: 16 -> 17

$ws = $wb.Worksheets.Item("Near South Side")
$ws.Cells.Item(5, 5).Value = 15  # E5: 14 -> 15
$ws.Cells.Item(6, 5).Value = 19  # E6: 18 -> 19

$ws = $wb.Worksheets.Item("Lower West Side")
$ws.Cells.Item(4, 8).Value = 2  # H4: 1 -> 2
$ws.Cells.Item(5, 3).Value = 2  # C5: 1 -> 2
$ws.Cells.Item(6, 3).Value = 6  # C6: 5 -> 6
$ws.Cells.Item(6, 8).Value = 4  # H6: 3 -> 4

$ws = $wb.Worksheets.Item("Roseland")
$ws.Cells.Item(7, 4).Value = 5  # D7: 4 -> 5
$ws.Cells.Item(8, 9).Value = 25  # I8: 24 -> 25
$ws.Cells.Item(9, 4).Value = 35  # D9: 34 -> 35
$ws.Cells.Item(9, 9).Value = 54  # I9: 53 -> 54

$ws = $wb.Worksheets.Item("Lincoln Square")
$ws.Cells.Item(5, 8).Value = 4  # H5: 3 -> 4
$ws.Cells.Item(6, 8).Value = 5  # H6: 4 -> 5

$ws = $wb.Worksheets.Item("Streeterville")
$ws.Cells.Item(5, 4).Value = 2  # D5: 1 -> 2
$ws.Cells.Item(6, 4).Value = 3  # D6: 2 -> 3

$ws = $wb.Worksheets.Item("Wicker Park")
$ws.Cells.Item(5, 6).Value = 11  # F5: 10 -> 11
$ws.Cells.Item(6, 6).Value = 13  # F6: 12 -> 13

$ws = $wb.Worksheets.Item("New City")
$ws.Cells.Item(5, 5).Value = 6  # E5: 5 -> 6
$ws.Cells.Item(6, 5).Value = 9  # E6: 8 -> 9

$ws = $wb.Worksheets.Item("Lincoln Park")
$ws.Cells.Item(6, 5).Value = 5  # E6: 4 -> 5
$ws.Cells.Item(7, 5).Value = 9  # E7: 8 -> 9

$ws = $wb.Worksheets.Item("Printers Row")
$ws.Cells.Item(3, 8).Value = 1  # H3: None -> 1
$ws.Cells.Item(6, 8).Value = 4  # H6: 3 -> 4

$ws = $wb.Worksheets.Item("Gage Park")
$ws.Cells.Item(3, 5).Value = 3  # E3: 2 -> 3
$ws.Cells.Item(6, 5).Value = 26  # E6: 25 -> 26

$ws = $wb.Worksheets.Item("Austin")
$ws.Cells.Item(7, 7).Value = 50  # G7: 49 -> 50
$ws.Cells.Item(8, 7).Value = 71  # G8: 70 -> 71

$ws = $wb.Worksheets.Item("Garfield Ridge")
$ws.Cells.Item(3, 2).Value = 2  # B3: 1 -> 2
$ws.Cells.Item(4, 2).Value = 2  # B4: 1 -> 2
